$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2748.3333
$ws.Range("I33").Value = 2875.25
$ws.Range("K33").Value = 2875.25
$ws.Range("M33").Value = -2646.25
$ws.Range("H135").Value = 17247440
$ws.Range("I135").Value = 25001228
$ws.Range("J135").Value = 16802.223
$ws.Range("K135").Value = 225011052
$ws.Range("L135").Value = 151220.007
$ws.Range("M135").Value = -225008517
$ws.Range("N135").Value = -156290.007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""
$ws.Range("H45").Value = 5235.2144
$ws.Range("I45").Value = 3897.8572
$ws.Range("J45").Value = 6572.5713
$ws.Range("K45").Value = 3897.8572
$ws.Range("L45").Value = 6572.5713
$ws.Range("M45").Value = -3520.8572
$ws.Range("N45").Value = -7326.5713
$ws.Range("H74").Value = 1177.1351
$ws.Range("I74").Value = 1155.8572
$ws.Range("K74").Value = 1155.8572
$ws.Range("M74").Value = -281.8571999999999
$ws.Range("H77").Value = 1177.1351
$ws.Range("I77").Value = 1155.8572
$ws.Range("K77").Value = 5779.286
$ws.Range("M77").Value = -1411.286
$ws.Range("H88").Value = 1770.9333
$ws.Range("I88").Value = 1797.4
$ws.Range("J88").Value = 1757.7
$ws.Range("K88").Value = 1797.4
$ws.Range("L88").Value = 1757.7
$ws.Range("M88").Value = -1391.4
$ws.Range("N88").Value = -2569.7
$ws.Range("H91").Value = 1770.9333
$ws.Range("I91").Value = 1797.4
$ws.Range("J91").Value = 1757.7
$ws.Range("K91").Value = 1797.4
$ws.Range("L91").Value = 1757.7
$ws.Range("M91").Value = -393.4000000000001
$ws.Range("N91").Value = -4565.7
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492
$ws.Range("H97").Value = 1239.6364
$ws.Range("I97").Value = 1322.5
$ws.Range("K97").Value = 1322.5
$ws.Range("M97").Value = -826.5
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H110").Value = 2879.5293
$ws.Range("I110").Value = 2871.9375
$ws.Range("J110").Value = 3001
$ws.Range("K110").Value = 2871.9375
$ws.Range("L110").Value = 3001
$ws.Range("M110").Value = -826.9375
$ws.Range("N110").Value = -7091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 58197.5
$ws.Range("J92").Value = 58197.5
$ws.Range("L92").Value = 58197.5
$ws.Range("N92").Value = -63189.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H103").Value = 21999.666
$ws.Range("J103").Value = 21999.666
$ws.Range("L103").Value = 21999.666
$ws.Range("N103").Value = -24343.666
$ws.Range("H104").Value = 90000
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H134").Value = 4147.7827
$ws.Range("I134").Value = 2037.3077
$ws.Range("K134").Value = 6111.9231
$ws.Range("M134").Value = -3576.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 159.66667
$ws.Range("I7").Value = 18.75
$ws.Range("K7").Value = 18.75
$ws.Range("M7").Value = 94.25
$ws.Range("H16").Value = 1122.7142
$ws.Range("I16").Value = 831.8
$ws.Range("J16").Value = 1850
$ws.Range("K16").Value = 831.8
$ws.Range("L16").Value = 1850
$ws.Range("M16").Value = -544.8
$ws.Range("N16").Value = -2424
$ws.Range("H22").Value = 520.7143
$ws.Range("I22").Value = 520.7143
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 520.7143
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -170.7143
$ws.Range("N22").Value = ""
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""
$ws.Range("H99").Value = 4320782.5
$ws.Range("I99").Value = 1631834.2
$ws.Range("J99").Value = 5410896.5
$ws.Range("K99").Value = 1631834.2
$ws.Range("L99").Value = 5410896.5
$ws.Range("M99").Value = -1630336.2
$ws.Range("N99").Value = -5413892.5
$ws.Range("H105").Value = 754.1875
$ws.Range("I105").Value = 676.8182
$ws.Range("K105").Value = 676.8182
$ws.Range("M105").Value = 1070.1818
$ws.Range("H113").Value = 1122.7142
$ws.Range("I113").Value = 831.8
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 831.8
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 1338.2
$ws.Range("N113").Value = -6190
$ws.Range("H122").Value = 301833
$ws.Range("I122").Value = 341659.1
$ws.Range("J122").Value = 3137.25
$ws.Range("K122").Value = 1024977.3
$ws.Range("L122").Value = 9411.75
$ws.Range("M122").Value = -1022527.3
$ws.Range("N122").Value = -14311.75
$ws.Range("H126").Value = 4320782.5
$ws.Range("I126").Value = 1631834.2
$ws.Range("J126").Value = 5410896.5
$ws.Range("K126").Value = 4895502.6
$ws.Range("L126").Value = 16232689.5
$ws.Range("M126").Value = -4893032.6
$ws.Range("N126").Value = -16237629.5
$ws.Range("H134").Value = 3698.0356
$ws.Range("I134").Value = 1850.4
$ws.Range("J134").Value = 5829.923
$ws.Range("K134").Value = 5551.200000000001
$ws.Range("L134").Value = 17489.769
$ws.Range("M134").Value = -3016.200000000001
$ws.Range("N134").Value = -22559.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 4997.8
$ws.Range("I103").Value = 4996.3335
$ws.Range("K103").Value = 14989.0005
$ws.Range("M103").Value = -14110.0005
$ws.Range("H132").Value = 612.3
$ws.Range("I132").Value = 580.7143
$ws.Range("J132").Value = 686
$ws.Range("K132").Value = 5226.428699999999
$ws.Range("L132").Value = 6174
$ws.Range("M132").Value = -2696.428699999999
$ws.Range("N132").Value = -11234

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5882450.5
$ws.Range("I2").Value = 87.40000000000001
$ws.Range("K2").Value = 87.40000000000001
$ws.Range("M2").Value = 25.59999999999999
$ws.Range("H39").Value = 40000
$ws.Range("J39").Value = 40000
$ws.Range("L39").Value = 40000
$ws.Range("N39").Value = -41064
$ws.Range("H80").Value = 53386.773
$ws.Range("I80").Value = 72097.625
$ws.Range("K80").Value = 72097.625
$ws.Range("M80").Value = -71099.625
$ws.Range("H83").Value = 53386.773
$ws.Range("I83").Value = 72097.625
$ws.Range("K83").Value = 360488.125
$ws.Range("M83").Value = -355496.125
$ws.Range("H102").Value = 2427.88
$ws.Range("I102").Value = 1795.238
$ws.Range("J102").Value = 5749.25
$ws.Range("K102").Value = 1795.238
$ws.Range("L102").Value = 5749.25
$ws.Range("M102").Value = -173.2380000000001
$ws.Range("N102").Value = -8993.25
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H122").Value = 4513.227
$ws.Range("J122").Value = 4433.3335
$ws.Range("L122").Value = 13300.0005
$ws.Range("N122").Value = -18200.0005
$ws.Range("H132").Value = 3003.5715
$ws.Range("I132").Value = 2669
$ws.Range("J132").Value = 5011
$ws.Range("K132").Value = 8007
$ws.Range("L132").Value = 15033
$ws.Range("M132").Value = -5477
$ws.Range("N132").Value = -20093

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1162.5
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505
$ws.Range("H27").Value = 1162.5
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693
$ws.Range("H30").Value = 2684.2
$ws.Range("I30").Value = 855.5
$ws.Range("J30").Value = 9999
$ws.Range("K30").Value = 855.5
$ws.Range("L30").Value = 9999
$ws.Range("M30").Value = -747.5
$ws.Range("N30").Value = -10215
$ws.Range("H32").Value = 15172.5
$ws.Range("I32").Value = 345
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 345
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -28
$ws.Range("N32").Value = -30634
$ws.Range("H35").Value = 1654.25
$ws.Range("I35").Value = 1654.25
$ws.Range("K35").Value = 1654.25
$ws.Range("M35").Value = -1318.25
$ws.Range("H46").Value = 8956.091
$ws.Range("I46").Value = 3309.8
$ws.Range("K46").Value = 3309.8
$ws.Range("M46").Value = -3121.8
$ws.Range("H136").Value = 3994.4187
$ws.Range("I136").Value = 2204.7273
$ws.Range("J136").Value = 5869.3335
$ws.Range("K136").Value = 6614.1819
$ws.Range("L136").Value = 17608.0005
$ws.Range("M136").Value = -4064.1819
$ws.Range("N136").Value = -22708.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4881.7334
$ws.Range("I122").Value = 2040.7727
$ws.Range("J122").Value = 12694.375
$ws.Range("K122").Value = 6122.3181
$ws.Range("L122").Value = 38083.125
$ws.Range("M122").Value = -3672.3181
$ws.Range("N122").Value = -42983.125
